# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New "Periodo Mora" (column E) values for rows 16-28 - now ascending order
$periodos = @("2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101")

# New "Valor Mora" (column F) values for rows 16-28 - rotated so the lower
# value (25013) now belongs to the last period (2101) instead of the first
$valorMora = @(31266,31266,31266,31266,31266,31266,31266,31266,31266,31266,31266,31266,25013)

# New "Salario Basico" (column G) values for rows 16-28 - all bumped by 400
$salario = @(781642,781642,781642,781642,781642,781642,781642,781642,781642,781642,781642,781642,781642)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valorMora[$i]
    $ws.Cells.Item($row, 7).Value = $salario[$i]
}

$wb.Save()
